$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Coefficient" column (E) holds numeric-looking text (e.g. "1", "2", "6")
# that must be stored as text, not as a number. Pre-formatting the range as
# Text before writing keeps Excel from auto-converting the values, and then
# resetting the style back to Normal restores the original (default) cell
# formatting/appearance.
$ws.Range("E3:E13").NumberFormat = "@"

# Row 3 - KIQ005 (How Often Do You Have Urinary Leakage per week?)
$ws.Range("D3").Value = " Ever day and/or night"
$ws.Range("E3").Value = "None"

# Row 4 - AUQ054 (General Condition of Hearing?)
$ws.Range("D4").Value = "Excellent"
$ws.Range("E4").Value = "1"

# Row 5 - MCQ560 (Have You Ever Had GallBladder Surgery?)
$ws.Range("D5").Value = "No"
$ws.Range("E5").Value = "2"

# Row 6 - MCQ371D (Are you watching your weight?)
$ws.Range("D6").Value = "No"
$ws.Range("E6").Value = "2"

# Row 8 - OHQ033 (Main Reason for Visiting Dentist?)
$ws.Range("D8").Value = "Went in on own for check-up, examination, or cleaning"
$ws.Range("E8").Value = "1"

# Row 10 - SMQ020 (Have You Smoked Atleast 100 Cigarettes?)
$ws.Range("D10").Value = "No"
$ws.Range("E10").Value = "2"

# Row 11 - DPQ040 (Over the Last Two Weeks have You Felt Tired or Had Little Energy?)
$ws.Range("D11").Value = "More than half the days"
$ws.Range("E11").Value = "2"

# Row 13 - HUQ051 (Number of Times Received Healthcare Over Past Year?)
$ws.Range("D13").Value = "10 to 12"
$ws.Range("E13").Value = "6"

# Restore the default cell style now that the text values are committed.
$ws.Range("E3:E13").Style = "Normal"
